$d = $word.ActiveDocument
$d.Content.Find.Execute("90+9=99", $true, $false, $false, $false, $false, $true, 1, $false, "77-70=7", 2) | Out-Null
$d.Content.Find.Execute("25+29=54", $true, $false, $false, $false, $false, $true, 1, $false, "64+9=73", 2) | Out-Null
$d.Content.Find.Execute("20+45=65", $true, $false, $false, $false, $false, $true, 1, $false, "44+18=62", 2) | Out-Null
$d.Content.Find.Execute("0+41=41", $true, $false, $false, $false, $false, $true, 1, $false, "31-2=29", 2) | Out-Null
$d.Content.Find.Execute("5+8=13", $true, $false, $false, $false, $false, $true, 1, $false, "95-90=5", 2) | Out-Null
$d.Content.Find.Execute("96-84=12", $true, $false, $false, $false, $false, $true, 1, $false, "64-5=59", 2) | Out-Null
$d.Content.Find.Execute("42+37=79", $true, $false, $false, $false, $false, $true, 1, $false, "18+64=82", 2) | Out-Null
$d.Content.Find.Execute("68+10=78", $true, $false, $false, $false, $false, $true, 1, $false, "49+46=95", 2) | Out-Null
$d.Content.Find.Execute("59-21=38", $true, $false, $false, $false, $false, $true, 1, $false, "12+79=91", 2) | Out-Null
$d.Content.Find.Execute("36+37=73", $true, $false, $false, $false, $false, $true, 1, $false, "44+34=78", 2) | Out-Null
$d.Content.Find.Execute("37+6=43", $true, $false, $false, $false, $false, $true, 1, $false, "81-50=31", 2) | Out-Null
$d.Content.Find.Execute("50+24=74", $true, $false, $false, $false, $false, $true, 1, $false, "60-59=1", 2) | Out-Null
$d.Content.Find.Execute("30+43=73", $true, $false, $false, $false, $false, $true, 1, $false, "61+33=94", 2) | Out-Null
$d.Content.Find.Execute("30+9=39", $true, $false, $false, $false, $false, $true, 1, $false, "42-15=27", 2) | Out-Null
$d.Content.Find.Execute("7+59=66", $true, $false, $false, $false, $false, $true, 1, $false, "57-2=55", 2) | Out-Null
$d.Content.Find.Execute("95-65=30", $true, $false, $false, $false, $false, $true, 1, $false, "68-42=26", 2) | Out-Null
$d.Content.Find.Execute("93-85=8", $true, $false, $false, $false, $false, $true, 1, $false, "75-28=47", 2) | Out-Null
$d.Content.Find.Execute("54-24=30", $true, $false, $false, $false, $false, $true, 1, $false, "15+43=58", 2) | Out-Null
$d.Content.Find.Execute("40+37=77", $true, $false, $false, $false, $false, $true, 1, $false, "19+53=72", 2) | Out-Null
$d.Content.Find.Execute("61+13=74", $true, $false, $false, $false, $false, $true, 1, $false, "78+20=98", 2) | Out-Null
$d.Content.Find.Execute("43-11=32", $true, $false, $false, $false, $false, $true, 1, $false, "8+16=24", 2) | Out-Null
$d.Content.Find.Execute("15+8=23", $true, $false, $false, $false, $false, $true, 1, $false, "84-31=53", 2) | Out-Null
$d.Content.Find.Execute("20+19=39", $true, $false, $false, $false, $false, $true, 1, $false, "21+13=34", 2) | Out-Null
$d.Content.Find.Execute("79-21=58", $true, $false, $false, $false, $false, $true, 1, $false, "45+41=86", 2) | Out-Null
$d.Content.Find.Execute("17-2=15", $true, $false, $false, $false, $false, $true, 1, $false, "13+8=21", 2) | Out-Null
$d.Content.Find.Execute("75-5=70", $true, $false, $false, $false, $false, $true, 1, $false, "44+52=96", 2) | Out-Null
$d.Content.Find.Execute("63-41=22", $true, $false, $false, $false, $false, $true, 1, $false, "94-63=31", 2) | Out-Null
$d.Content.Find.Execute("94-29=65", $true, $false, $false, $false, $false, $true, 1, $false, "10-8=2", 2) | Out-Null
$d.Content.Find.Execute("61+36=97", $true, $false, $false, $false, $false, $true, 1, $false, "98-2=96", 2) | Out-Null
$d.Content.Find.Execute("21-4=17", $true, $false, $false, $false, $false, $true, 1, $false, "78-11=67", 2) | Out-Null
$d.Content.Find.Execute("53+2=55", $true, $false, $false, $false, $false, $true, 1, $false, "10+72=82", 2) | Out-Null
$d.Content.Find.Execute("79-0=79", $true, $false, $false, $false, $false, $true, 1, $false, "49-24=25", 2) | Out-Null
$d.Content.Find.Execute("11+50=61", $true, $false, $false, $false, $false, $true, 1, $false, "33-30=3", 2) | Out-Null
$d.Content.Find.Execute("56-19=37", $true, $false, $false, $false, $false, $true, 1, $false, "68-27=41", 2) | Out-Null
$d.Content.Find.Execute("95-88=7", $true, $false, $false, $false, $false, $true, 1, $false, "59-0=59", 2) | Out-Null
$d.Content.Find.Execute("81-35=46", $true, $false, $false, $false, $false, $true, 1, $false, "32+60=92", 2) | Out-Null
$d.Content.Find.Execute("58-26=32", $true, $false, $false, $false, $false, $true, 1, $false, "28+22=50", 2) | Out-Null
$d.Content.Find.Execute("50+33=83", $true, $false, $false, $false, $false, $true, 1, $false, "41-36=5", 2) | Out-Null
$d.Content.Find.Execute("1+92=93", $true, $false, $false, $false, $false, $true, 1, $false, "52-10=42", 2) | Out-Null
$d.Content.Find.Execute("85+0=85", $true, $false, $false, $false, $false, $true, 1, $false, "6+50=56", 2) | Out-Null
$d.Content.Find.Execute("11+40=51", $true, $false, $false, $false, $false, $true, 1, $false, "12+14=26", 2) | Out-Null
$d.Content.Find.Execute("55+34=89", $true, $false, $false, $false, $false, $true, 1, $false, "3+22=25", 2) | Out-Null
$d.Content.Find.Execute("43+37=80", $true, $false, $false, $false, $false, $true, 1, $false, "12+81=93", 2) | Out-Null
$d.Content.Find.Execute("17+9=26", $true, $false, $false, $false, $false, $true, 1, $false, "51+20=71", 2) | Out-Null
$d.Content.Find.Execute("10+8=18", $true, $false, $false, $false, $false, $true, 1, $false, "28+58=86", 2) | Out-Null
$d.Content.Find.Execute("20+13=33", $true, $false, $false, $false, $false, $true, 1, $false, "42-39=3", 2) | Out-Null
$d.Content.Find.Execute("20+40=60", $true, $false, $false, $false, $false, $true, 1, $false, "40-14=26", 2) | Out-Null
$d.Content.Find.Execute("51-13=38", $true, $false, $false, $false, $false, $true, 1, $false, "95-33=62", 2) | Out-Null
$d.Content.Find.Execute("60-30=30", $true, $false, $false, $false, $false, $true, 1, $false, "81-43=38", 2) | Out-Null
$d.Content.Find.Execute("13-0=13", $true, $false, $false, $false, $false, $true, 1, $false, "91-43=48", 2) | Out-Null
$d.Content.Find.Execute("70-7=63", $true, $false, $false, $false, $false, $true, 1, $false, "27+22=49", 2) | Out-Null
$d.Content.Find.Execute("38+33=71", $true, $false, $false, $false, $false, $true, 1, $false, "20-0=20", 2) | Out-Null
$d.Content.Find.Execute("69-13=56", $true, $false, $false, $false, $false, $true, 1, $false, "6+69=75", 2) | Out-Null
$d.Content.Find.Execute("90-61=29", $true, $false, $false, $false, $false, $true, 1, $false, "36+36=72", 2) | Out-Null
$d.Content.Find.Execute("36+21=57", $true, $false, $false, $false, $false, $true, 1, $false, "90-43=47", 2) | Out-Null
$d.Content.Find.Execute("7+41=48", $true, $false, $false, $false, $false, $true, 1, $false, "73-9=64", 2) | Out-Null
$d.Content.Find.Execute("3+93=96", $true, $false, $false, $false, $false, $true, 1, $false, "17+7=24", 2) | Out-Null
$d.Content.Find.Execute("90-41=49", $true, $false, $false, $false, $false, $true, 1, $false, "88-22=66", 2) | Out-Null
$d.Content.Find.Execute("75+6=81", $true, $false, $false, $false, $false, $true, 1, $false, "87-1=86", 2) | Out-Null
$d.Content.Find.Execute("2+80=82", $true, $false, $false, $false, $false, $true, 1, $false, "35+15=50", 2) | Out-Null
$d.Content.Find.Execute("76-75=1", $true, $false, $false, $false, $false, $true, 1, $false, "67-38=29", 2) | Out-Null
$d.Content.Find.Execute("76-59=17", $true, $false, $false, $false, $false, $true, 1, $false, "50+39=89", 2) | Out-Null
$d.Content.Find.Execute("80-45=35", $true, $false, $false, $false, $false, $true, 1, $false, "88-35=53", 2) | Out-Null
$d.Content.Find.Execute("13+18=31", $true, $false, $false, $false, $false, $true, 1, $false, "71-1=70", 2) | Out-Null
$d.Content.Find.Execute("42+15=57", $true, $false, $false, $false, $false, $true, 1, $false, "25+55=80", 2) | Out-Null
$d.Content.Find.Execute("53-14=39", $true, $false, $false, $false, $false, $true, 1, $false, "88-14=74", 2) | Out-Null
$d.Content.Find.Execute("92-86=6", $true, $false, $false, $false, $false, $true, 1, $false, "88+4=92", 2) | Out-Null
$d.Content.Find.Execute("97-61=36", $true, $false, $false, $false, $false, $true, 1, $false, "3+15=18", 2) | Out-Null
$d.Content.Find.Execute("40+35=75", $true, $false, $false, $false, $false, $true, 1, $false, "65-10=55", 2) | Out-Null
$d.Content.Find.Execute("28+32=60", $true, $false, $false, $false, $false, $true, 1, $false, "69-10=59", 2) | Out-Null
$d.Content.Find.Execute("47+44=91", $true, $false, $false, $false, $false, $true, 1, $false, "96-94=2", 2) | Out-Null
$d.Content.Find.Execute("46-37=9", $true, $false, $false, $false, $false, $true, 1, $false, "77-61=16", 2) | Out-Null
$d.Content.Find.Execute("62-18=44", $true, $false, $false, $false, $false, $true, 1, $false, "32+16=48", 2) | Out-Null
$d.Content.Find.Execute("40-39=1", $true, $false, $false, $false, $false, $true, 1, $false, "3+0=3", 2) | Out-Null
$d.Content.Find.Execute("31+42=73", $true, $false, $false, $false, $false, $true, 1, $false, "42-36=6", 2) | Out-Null
$d.Content.Find.Execute("65-42=23", $true, $false, $false, $false, $false, $true, 1, $false, "75-46=29", 2) | Out-Null
$d.Content.Find.Execute("5+79=84", $true, $false, $false, $false, $false, $true, 1, $false, "90-22=68", 2) | Out-Null
$d.Content.Find.Execute("11+18=29", $true, $false, $false, $false, $false, $true, 1, $false, "70+0=70", 2) | Out-Null
$d.Content.Find.Execute("30+45=75", $true, $false, $false, $false, $false, $true, 1, $false, "11+53=64", 2) | Out-Null
$d.Content.Find.Execute("47+49=96", $true, $false, $false, $false, $false, $true, 1, $false, "93-15=78", 2) | Out-Null
$d.Content.Find.Execute("68-44=24", $true, $false, $false, $false, $false, $true, 1, $false, "61-11=50", 2) | Out-Null
$d.Content.Find.Execute("46+32=78", $true, $false, $false, $false, $false, $true, 1, $false, "42+33=75", 2) | Out-Null
$d.Content.Find.Execute("65-15=50", $true, $false, $false, $false, $false, $true, 1, $false, "60-4=56", 2) | Out-Null
$d.Content.Find.Execute("73-53=20", $true, $false, $false, $false, $false, $true, 1, $false, "22+17=39", 2) | Out-Null
$d.Content.Find.Execute("38-14=24", $true, $false, $false, $false, $false, $true, 1, $false, "0+91=91", 2) | Out-Null
$d.Content.Find.Execute("23+46=69", $true, $false, $false, $false, $false, $true, 1, $false, "96-83=13", 2) | Out-Null
$d.Content.Find.Execute("21+24=45", $true, $false, $false, $false, $false, $true, 1, $false, "10+22=32", 2) | Out-Null
$d.Content.Find.Execute("26+52=78", $true, $false, $false, $false, $false, $true, 1, $false, "39+1=40", 2) | Out-Null
$d.Content.Find.Execute("19+63=82", $true, $false, $false, $false, $false, $true, 1, $false, "54-28=26", 2) | Out-Null
$d.Content.Find.Execute("31+43=74", $true, $false, $false, $false, $false, $true, 1, $false, "59-51=8", 2) | Out-Null
$d.Content.Find.Execute("50+48=98", $true, $false, $false, $false, $false, $true, 1, $false, "18+9=27", 2) | Out-Null
$d.Content.Find.Execute("69-24=45", $true, $false, $false, $false, $false, $true, 1, $false, "73+19=92", 2) | Out-Null
$d.Content.Find.Execute("30-1=29", $true, $false, $false, $false, $false, $true, 1, $false, "10+69=79", 2) | Out-Null
$d.Content.Find.Execute("17+22=39", $true, $false, $false, $false, $false, $true, 1, $false, "68-57=11", 2) | Out-Null
$d.Content.Find.Execute("49+18=67", $true, $false, $false, $false, $false, $true, 1, $false, "96-52=44", 2) | Out-Null
$d.Content.Find.Execute("7+15=22", $true, $false, $false, $false, $false, $true, 1, $false, "88-71=17", 2) | Out-Null
$d.Content.Find.Execute("23+61=84", $true, $false, $false, $false, $false, $true, 1, $false, "98-35=63", 2) | Out-Null
$d.Content.Find.Execute("44-3=41", $true, $false, $false, $false, $false, $true, 1, $false, "25-19=6", 2) | Out-Null
$d.Content.Find.Execute("75+21=96", $true, $false, $false, $false, $false, $true, 1, $false, "5+31=36", 2) | Out-Null
$d.Content.Find.Execute("62-20=42", $true, $false, $false, $false, $false, $true, 1, $false, "91-59=32", 2) | Out-Null
